$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.878.92"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").Value = "2.632.37"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.75"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.62"
$ws.Range("E6").Value = "  +2.15%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +5.17%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.398"
$ws.Range("E10").Value = "  +3.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  +1.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("E12").Value = "  +1.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.92"
$ws.Range("E13").Value = "  +3.35%  "

# Row 14
$ws.Range("D14").Value = "3.101.42"
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +12.54%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.711.64"
$ws.Range("E16").Value = "  +1.61%  "

# Row 17
$ws.Range("D17").Value = "2.623.49"
$ws.Range("E17").Value = "  -1.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.50"
$ws.Range("E18").Value = "  +0.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  +2.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.18"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.19"
$ws.Range("E21").Value = "  +4.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.70"
$ws.Range("E23").Value = "  +0.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.44"
$ws.Range("E25").Value = "  +1.03%  "

# Row 26
$ws.Range("E26").Value = "  -1.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.33"
$ws.Range("E27").Value = "  +2.16%  "

# Row 28
$ws.Range("E28").Value = "  +1.25%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.40%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "539.16"
$ws.Range("E30").Value = "  -3.20%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0919"
$ws.Range("E31").Value = "  +7.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  +0.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("E33").Value = "  +3.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  +7.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.27"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").Value = "  +1.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.93"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38
$ws.Range("E38").Value = "  +3.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.13"
$ws.Range("E39").Value = "  +2.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.78"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.92"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.11"
$ws.Range("E44").Value = "  +3.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.52"
$ws.Range("E45").Value = "  +6.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0603"
$ws.Range("E46").Value = "  +2.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("E47").Value = "  +9.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.642"
$ws.Range("E48").Value = "  +1.18%  "

# Row 49
$ws.Range("E49").Value = "  -0.83%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0984"
$ws.Range("E50").Value = "  +1.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.40"
$ws.Range("E51").Value = "  -0.49%  "
